$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (and workbook title reflects latest date through 10-27)
$ws.Name = "Through 2021-10-27"

# Row 12 - October (through 10-26) -> October (through 10-27) with updated counts
$ws.Range("A12").Value = "October (through 10-27)"
$ws.Range("C12").Value = 25
$ws.Range("D12").Value = 0.0741
$ws.Range("I12").Value = 59
$ws.Range("J12").Value = 0.1194
$ws.Range("L12").Value = 52
$ws.Range("M12").Value = 0.0877
$ws.Range("O12").Value = 48
$ws.Range("P12").Value = 0.0769

# New cells Q12 and S12 need to match the percentage style used by column D/G/J/M/P/S/V (s="1", format 0.0%)
$ws.Range("Q12").NumberFormat = "0.0%"
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = 132
$ws.Range("S12").NumberFormat = "0.0%"
$ws.Range("S12").Value = 0.0075

$ws.Range("U12").Value = 173

# Row 13 - Total row updates
$ws.Range("C13").Value = 221
$ws.Range("D13").Value = 0.1265
$ws.Range("I13").Value = 636
$ws.Range("J13").Value = 0.08359999999999999
$ws.Range("L13").Value = 539
$ws.Range("M13").Value = 0.1091
$ws.Range("O13").Value = 427
$ws.Range("P13").Value = 0.0992
$ws.Range("Q13").Value = 54
$ws.Range("R13").Value = 980
$ws.Range("S13").Value = 0.0522
$ws.Range("U13").Value = 1338
$ws.Range("V13").Value = 0.0577
